$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter data the way the author originally typed it: down column A first,
# then the header for column B, then the last row of column A.
$ws.Range("A1").Value = "Things That Need Work"
$ws.Range("A2").Value = "The images/styles need to work"
$ws.Range("B1").Value = "Status"
$ws.Range("A3").Value = "DB needs to be modified in order to accept login screen parameters"

# Column width to fit the long text in column A (best-fit of the longest
# entry, "DB needs to be modified in order to accept login screen parameters")
$ws.Columns.Item(1).ColumnWidth = 59.6

# Turn the range into a real Excel Table (ListObject)
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:B3"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"

# Selection as left by the author after entering data
$ws.Range("A4").Select() | Out-Null
